# Add season record columns (Wins, Losses, Ties) to the roster/stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels in AD1:AF1, matching the header style used
#     by the rest of row 1 (bold, bordered, centered) by copying AC1's format.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-52): every player row gets the same season record values.
$wins = 79
$losses = 83
$ties = 0

for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-52"
